$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 178, shifting existing rows 178-252 down to 179-253.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with the new weekly record.
$ws.Cells.Item(178, 1).Value = 9
$ws.Cells.Item(178, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(178, 3).Value = "Metropolitana"
$ws.Cells.Item(178, 4).Value = 44466
$ws.Cells.Item(178, 5).Value = 13
$ws.Cells.Item(178, 6).Value = 100112031
$ws.Cells.Item(178, 7).Value = "Poroto verde"
$ws.Cells.Item(178, 8).Value = "Magnum"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 16
$ws.Cells.Item(178, 11).Value = 38000
$ws.Cells.Item(178, 12).Value = 40000
$ws.Cells.Item(178, 13).Value = 39000
$ws.Cells.Item(178, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(178, 15).Value = "Perú"
$ws.Cells.Item(178, 16).Value = 1560
$ws.Cells.Item(178, 17).Value = 25
$ws.Cells.Item(178, 18).Value = "Hortaliza"
